$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handback" - the localization-status report is
# refreshed: the zh-cn and de-de rows now show a completed handback (target
# file link + handback file name + handback datetime), the overall status
# flips from "Ready for handoff" to "Handed back: in sync with en-US", and a
# couple of report columns are widened to fit the longer text.
# ---------------------------------------------------------------------------

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status column: every cell that said "Ready for handoff" now reads the
# new handback status (Overview!E:F, and the "Status" column on each locale
# sheet). ---------------------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Per-locale handback details --------------------------------------
$base = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c8b1c1dbdb4743181aedb6c0de7b2069f45e86ca/e2e/"

$zhHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

# zh-cn: Latest Target File (I) becomes a link to a.md, Latest Handback File
# (J) is filled in, Latest Handback DateTime (K) gets stamped.
$wsZhCn.Range("J2").Value = $zhHandbackFile
$wsZhCn.Range("J3").Value = $zhHandbackFile
$wsZhCn.Range("K2").Value = "2016-08-30 04:36:45"
$wsZhCn.Range("K3").Value = "2016-08-30 04:36:45"

# de-de: same shape, different file suffix + handback timestamp.
$wsDeDe.Range("J2").Value = $deHandbackFile
$wsDeDe.Range("J3").Value = $deHandbackFile
$wsDeDe.Range("K2").Value = "2016-08-30 04:36:52"
$wsDeDe.Range("K3").Value = "2016-08-30 04:36:52"

# Rebuild hyperlinks so the new "Latest Target File" links land between the
# existing "Source File Name" links, in row order (A2, I2, A3, I3).
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), ($base + "a.md"), $null, $null, "a.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), ($base + "a.md"), $null, $null, "a.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), ($base + "b.md"), $null, $null, "b.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), ($base + "a.md"), $null, $null, "a.md") | Out-Null

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), ($base + "a.md"), $null, $null, "a.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), ($base + "a.md"), $null, $null, "a.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), ($base + "b.md"), $null, $null, "b.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), ($base + "a.md"), $null, $null, "a.md") | Out-Null

# --- Widen columns that now hold longer text ----------------------------
# (raw xlsx column width = ColumnWidth + 5/6; pick ColumnWidth values that
# land on the target widths of ~30 and 40 characters respectively.)
$wideWidth = 29.166666666666668   # -> stored width 30
$fortyWidth = 39.166666666666664  # -> stored width 40

$wsOverview.Columns.Item("E").ColumnWidth = $wideWidth
$wsOverview.Columns.Item("F").ColumnWidth = $wideWidth

$wsZhCn.Columns.Item("C").ColumnWidth = $wideWidth
$wsZhCn.Columns.Item("J").ColumnWidth = $fortyWidth

$wsDeDe.Columns.Item("C").ColumnWidth = $wideWidth
$wsDeDe.Columns.Item("J").ColumnWidth = $fortyWidth
